$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.5487310763284654
$ws.Cells.Item(2, 3).Value = 0.04628598318475952
$ws.Cells.Item(2, 4).Value = 0.06813507012373599
$ws.Cells.Item(2, 6).Value = 1.68530142800185
$ws.Cells.Item(2, 7).Value = 0.002498868747215143
$ws.Cells.Item(2, 11).Value = 0.5177818012245723
$ws.Cells.Item(2, 13).Value = 0.8774898827038982
$ws.Cells.Item(2, 14).Value = 2.554056340231384
$ws.Cells.Item(3, 2).Value = 0.5120698686170897
$ws.Cells.Item(3, 3).Value = 0.04092766014900917
$ws.Cells.Item(3, 4).Value = 0.06827354830565069
$ws.Cells.Item(3, 6).Value = 1.65590589052303
$ws.Cells.Item(3, 7).Value = 0.002502912754646737
$ws.Cells.Item(3, 11).Value = 0.4783511471295867
$ws.Cells.Item(3, 13).Value = 0.7818136211122635
$ws.Cells.Item(3, 14).Value = 2.549120797849866
$ws.Cells.Item(4, 2).Value = 0.489886788269132
$ws.Cells.Item(4, 3).Value = 0.03765358004258701
$ws.Cells.Item(4, 4).Value = 0.06835520188814881
$ws.Cells.Item(4, 6).Value = 1.638639528467422
$ws.Cells.Item(4, 7).Value = 0.002505524586723373
$ws.Cells.Item(4, 11).Value = 0.4544399296461563
$ws.Cells.Item(4, 13).Value = 0.7236161133314596
$ws.Cells.Item(4, 14).Value = 2.546719123599146
$ws.Cells.Item(5, 2).Value = 0.4809290774219903
$ws.Cells.Item(5, 3).Value = 0.03632331954131018
$ws.Cells.Item(5, 4).Value = 0.06838763115588997
$ws.Cells.Item(5, 6).Value = 1.63179948471155
$ws.Cells.Item(5, 7).Value = 0.002506621430187587
$ws.Cells.Item(5, 11).Value = 0.4447709615452595
$ws.Cells.Item(5, 13).Value = 0.700031484842242
$ws.Cells.Item(5, 14).Value = 2.545897926150275
$ws.Cells.Item(6, 2).Value = 0.4794466088263789
$ws.Cells.Item(6, 3).Value = 0.03610266793289441
$ws.Cells.Item(6, 4).Value = 0.06839296504527859
$ws.Cells.Item(6, 6).Value = 1.630675528004559
$ws.Cells.Item(6, 7).Value = 0.002506805526518616
$ws.Cells.Item(6, 11).Value = 0.4431699628792387
$ws.Cells.Item(6, 13).Value = 0.6961230195227586
$ws.Cells.Item(6, 14).Value = 2.545771062252413
$ws.Cells.Item(7, 2).Value = 0.4897656493419049
$ws.Cells.Item(7, 3).Value = 0.03763562375176832
$ws.Cells.Item(7, 4).Value = 0.06835564265911631
$ws.Cells.Item(7, 6).Value = 1.638546487940161
$ws.Cells.Item(7, 7).Value = 0.002505539247464312
$ws.Cells.Item(7, 11).Value = 0.4543092269954059
$ws.Cells.Item(7, 13).Value = 0.7232975195421432
$ws.Cells.Item(7, 14).Value = 2.546707411703082
$ws.Cells.Item(8, 2).Value = 0.5360222992001127
$ws.Cells.Item(8, 3).Value = 0.04443508612943958
$ws.Cells.Item(8, 4).Value = 0.06818351908998022
$ws.Cells.Item(8, 6).Value = 1.675003004320075
$ws.Cells.Item(8, 7).Value = 0.002500236458389892
$ws.Cells.Item(8, 11).Value = 0.504123825121269
$ws.Cells.Item(8, 13).Value = 0.8443829763395883
$ws.Cells.Item(8, 14).Value = 2.552223694005207
$ws.Cells.Item(9, 2).Value = 0.6293394883501264
$ws.Cells.Item(9, 3).Value = 0.05789868454419889
$ws.Cells.Item(9, 4).Value = 0.06781908261914804
$ws.Cells.Item(9, 6).Value = 1.75273915454261
$ws.Cells.Item(9, 7).Value = 0.002490854421654742
$ws.Cells.Item(9, 11).Value = 0.6042016576192282
$ws.Cells.Item(9, 13).Value = 1.086489041410047
$ws.Cells.Item(9, 14).Value = 2.56805955274514
$ws.Cells.Item(10, 2).Value = 0.6995157349705039
$ws.Cells.Item(10, 3).Value = 0.06787515685707035
$ws.Cells.Item(10, 4).Value = 0.06753471114104714
$ws.Cells.Item(10, 6).Value = 1.813716215835981
$ws.Cells.Item(10, 7).Value = 0.002484573929688946
$ws.Cells.Item(10, 11).Value = 0.6792177358785523
$ws.Cells.Item(10, 13).Value = 1.267655065277381
$ws.Cells.Item(10, 14).Value = 2.582797471587156
$ws.Cells.Item(11, 2).Value = 0.7317978812056936
$ws.Cells.Item(11, 3).Value = 0.07243339530008086
$ws.Cells.Item(11, 4).Value = 0.06740168847228656
$ws.Cells.Item(11, 6).Value = 1.842308381441597
$ws.Cells.Item(11, 7).Value = 0.002481848198873489
$ws.Cells.Item(11, 11).Value = 0.7136750691019245
$ws.Cells.Item(11, 13).Value = 1.350887973071139
$ws.Cells.Item(11, 14).Value = 2.590185812804009
$ws.Cells.Item(12, 2).Value = 0.7440742038728843
$ws.Cells.Item(12, 3).Value = 0.07416242439084897
$ws.Cells.Item(12, 4).Value = 0.06735078711021458
$ws.Cells.Item(12, 6).Value = 1.853259122818741
$ws.Cells.Item(12, 7).Value = 0.00248083479548877
$ws.Cells.Item(12, 11).Value = 0.7267713476909421
$ws.Cells.Item(12, 13).Value = 1.382532127746074
$ws.Cells.Item(12, 14).Value = 2.593082704718057
$ws.Cells.Item(13, 2).Value = 0.7414279720985348
$ws.Cells.Item(13, 3).Value = 0.07378991602280394
$ws.Cells.Item(13, 4).Value = 0.06736177316958081
$ws.Cells.Item(13, 6).Value = 1.850895178867205
$ws.Cells.Item(13, 7).Value = 0.002481052216948996
$ws.Cells.Item(13, 11).Value = 0.7239486890108822
$ws.Cells.Item(13, 13).Value = 1.375711277676729
$ws.Cells.Item(13, 14).Value = 2.592454388492115
$ws.Cells.Item(14, 2).Value = 0.7328068228179347
$ws.Cells.Item(14, 3).Value = 0.07257558461802205
$ws.Cells.Item(14, 4).Value = 0.06739751139180683
$ws.Cells.Item(14, 6).Value = 1.843206825636827
$ws.Cells.Item(14, 7).Value = 0.002481764449904405
$ws.Cells.Item(14, 11).Value = 0.7147515435544562
$ws.Cells.Item(14, 13).Value = 1.353488792884676
$ws.Cells.Item(14, 14).Value = 2.590422151758247
$ws.Cells.Item(15, 2).Value = 0.727532868749023
$ws.Cells.Item(15, 3).Value = 0.07183215412392485
$ws.Cells.Item(15, 4).Value = 0.06741933318838367
$ws.Cells.Item(15, 6).Value = 1.838513595532405
$ws.Cells.Item(15, 7).Value = 0.002482203154604345
$ws.Cells.Item(15, 11).Value = 0.7091242901104522
$ws.Cells.Item(15, 13).Value = 1.33989347895097
$ws.Cells.Item(15, 14).Value = 2.589190274080948
$ws.Cells.Item(16, 2).Value = 0.6974132545274188
$ws.Cells.Item(16, 3).Value = 0.06757767086116928
$ws.Cells.Item(16, 4).Value = 0.0675433306345079
$ws.Cells.Item(16, 6).Value = 1.811864890380406
$ws.Cells.Item(16, 7).Value = 0.00248475469538759
$ws.Cells.Item(16, 11).Value = 0.6769725748036421
$ws.Cells.Item(16, 13).Value = 1.262232753498296
$ws.Cells.Item(16, 14).Value = 2.582328451891129
$ws.Cells.Item(17, 2).Value = 0.679027841665544
$ws.Cells.Item(17, 3).Value = 0.06497282216410838
$ws.Cells.Item(17, 4).Value = 0.06761845987778736
$ws.Cells.Item(17, 6).Value = 1.79573590624554
$ws.Cells.Item(17, 7).Value = 0.002486353534041016
$ws.Cells.Item(17, 11).Value = 0.6573338001104219
$ws.Cells.Item(17, 13).Value = 1.214805815945397
$ws.Cells.Item(17, 14).Value = 2.578294648081979
$ws.Cells.Item(18, 2).Value = 0.6684867807000217
$ws.Cells.Item(18, 3).Value = 0.06347645052511552
$ws.Cells.Item(18, 4).Value = 0.06766132779288814
$ws.Cells.Item(18, 6).Value = 1.786539218467368
$ws.Cells.Item(18, 7).Value = 0.002487285508057369
$ws.Cells.Item(18, 11).Value = 0.646069342425875
$ws.Cells.Item(18, 13).Value = 1.187603759088077
$ws.Cells.Item(18, 14).Value = 2.576038846603851
$ws.Cells.Item(19, 2).Value = 0.6649235477015907
$ws.Cells.Item(19, 3).Value = 0.06297012380760236
$ws.Cells.Item(19, 4).Value = 0.06767578304333899
$ws.Cells.Item(19, 6).Value = 1.783439141823322
$ws.Cells.Item(19, 7).Value = 0.002487603185589199
$ws.Cells.Item(19, 11).Value = 0.6422607489733991
$ws.Cells.Item(19, 13).Value = 1.178406542542476
$ws.Cells.Item(19, 14).Value = 2.575286098626435
$ws.Cells.Item(20, 2).Value = 0.6809815065939233
$ws.Cells.Item(20, 3).Value = 0.06524991892673881
$ws.Cells.Item(20, 4).Value = 0.06761049790653573
$ws.Cells.Item(20, 6).Value = 1.797444550061002
$ws.Cells.Item(20, 7).Value = 0.002486182055848019
$ws.Cells.Item(20, 11).Value = 0.659421145401069
$ws.Cells.Item(20, 13).Value = 1.21984650416897
$ws.Cells.Item(20, 14).Value = 2.578717389622497
$ws.Cells.Item(21, 2).Value = 0.7353376565214091
$ws.Cells.Item(21, 3).Value = 0.07293218335007623
$ws.Cells.Item(21, 4).Value = 0.06738702857086665
$ws.Cells.Item(21, 6).Value = 1.845461722969347
$ws.Cells.Item(21, 7).Value = 0.002481554741204354
$ws.Cells.Item(21, 11).Value = 0.7174516622688429
$ws.Cells.Item(21, 13).Value = 1.360012601782728
$ws.Cells.Item(21, 14).Value = 2.591016373856945
$ws.Cells.Item(22, 2).Value = 0.7711644249128256
$ws.Cells.Item(22, 3).Value = 0.07797005524275846
$ws.Cells.Item(22, 4).Value = 0.06723789667125457
$ws.Cells.Item(22, 6).Value = 1.877564003427636
$ws.Cells.Item(22, 7).Value = 0.002478639889592316
$ws.Cells.Item(22, 11).Value = 0.7556580725693038
$ws.Cells.Item(22, 13).Value = 1.45235545602668
$ws.Cells.Item(22, 14).Value = 2.599632376411165
$ws.Cells.Item(23, 2).Value = 0.7520153122655984
$ws.Cells.Item(23, 3).Value = 0.07527966548130394
$ws.Cells.Item(23, 4).Value = 0.0673177738304247
$ws.Cells.Item(23, 6).Value = 1.860364240705337
$ws.Cells.Item(23, 7).Value = 0.002480185629615917
$ws.Cells.Item(23, 11).Value = 0.73524086988877
$ws.Cells.Item(23, 13).Value = 1.403000341111763
$ws.Cells.Item(23, 14).Value = 2.594980730618261
$ws.Cells.Item(24, 2).Value = 0.6800981644976503
$ws.Cells.Item(24, 3).Value = 0.06512463977644245
$ws.Cells.Item(24, 4).Value = 0.06761409852524825
$ws.Cells.Item(24, 6).Value = 1.796671835279241
$ws.Cells.Item(24, 7).Value = 0.002486259541190768
$ws.Cells.Item(24, 11).Value = 0.6584773752091166
$ws.Cells.Item(24, 13).Value = 1.217567409264433
$ws.Cells.Item(24, 14).Value = 2.578526071105117
$ws.Cells.Item(25, 2).Value = 0.6038128284637594
$ws.Cells.Item(25, 3).Value = 0.05424194591606124
$ws.Cells.Item(25, 4).Value = 0.06792057575338362
$ws.Cells.Item(25, 6).Value = 1.731035033182906
$ws.Cells.Item(25, 7).Value = 0.002493284423389768
$ws.Cells.Item(25, 11).Value = 0.5768687494065432
$ws.Cells.Item(25, 13).Value = 1.020448294029549
$ws.Cells.Item(25, 14).Value = 2.563233772840164
